# Update LR-pair metrics with newly recomputed TPM-based values.
# (commit: "update scripts wuth new tpm")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value  = 0.04575233333333333
$ws.Range("H2").Value  = 0.137257
$ws.Range("I2").Value  = 0.14147347546269
$ws.Range("J2").Value  = 0.14147347546269
$ws.Range("M2").Value  = 61.04160633333334
$ws.Range("N2").Value  = 183.124819
$ws.Range("O2").Value  = 0.2043613460574534
$ws.Range("P2").Value  = 0.2043613460574534
$ws.Range("Q2").Value  = 2.792795920164778
$ws.Range("R2").Value  = 25.135163281483
$ws.Range("S2").Value  = 0.02891170987698144
$ws.Range("T2").Value  = 0.02891170987698144

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value  = 0.04575233333333333
$ws.Range("H3").Value  = 0.137257
$ws.Range("I3").Value  = 0.14147347546269
$ws.Range("J3").Value  = 0.14147347546269
$ws.Range("O3").Value  = 0.3559304658284363
$ws.Range("P3").Value  = 0.3559304658284363
$ws.Range("Q3").Value  = 4.864134886587333
$ws.Range("R3").Value  = 43.777213979286
$ws.Range("S3").Value  = 0.05035472002380311
$ws.Range("T3").Value  = 0.0503547200238031

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value  = 0.04575233333333333
$ws.Range("H4").Value  = 0.137257
$ws.Range("I4").Value  = 0.14147347546269
$ws.Range("J4").Value  = 0.14147347546269
$ws.Range("M4").Value  = 131.3384093333333
$ws.Range("N4").Value  = 394.015228
$ws.Range("O4").Value  = 0.4397081881141102
$ws.Range("P4").Value  = 0.4397081881141103
$ws.Range("Q4").Value  = 6.009038683288444
$ws.Range("R4").Value  = 54.081348149596
$ws.Range("S4").Value  = 0.06220704556190546
$ws.Range("T4").Value  = 0.06220704556190546

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value  = 0.2599907647526892
$ws.Range("J5").Value  = 0.2599907647526892
$ws.Range("M5").Value  = 61.04160633333334
$ws.Range("N5").Value  = 183.124819
$ws.Range("O5").Value  = 0.2043613460574534
$ws.Range("P5").Value  = 0.2043613460574534
$ws.Range("Q5").Value  = 5.13241895491089
$ws.Range("R5").Value  = 46.19177059419801
$ws.Range("S5").Value  = 0.05313206264736628
$ws.Range("T5").Value  = 0.05313206264736628

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value  = 0.2599907647526892
$ws.Range("J6").Value  = 0.2599907647526892
$ws.Range("O6").Value  = 0.3559304658284363
$ws.Range("P6").Value  = 0.3559304658284363
$ws.Range("S6").Value  = 0.09253863400951605
$ws.Range("T6").Value  = 0.09253863400951606

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value  = 0.2599907647526892
$ws.Range("J7").Value  = 0.2599907647526892
$ws.Range("M7").Value  = 131.3384093333333
$ws.Range("N7").Value  = 394.015228
$ws.Range("O7").Value  = 0.4397081881141102
$ws.Range("P7").Value  = 0.4397081881141103
$ws.Range("Q7").Value  = 11.04302101568622
$ws.Range("R7").Value  = 99.38718914117601
$ws.Range("S7").Value  = 0.1143200680958068
$ws.Range("T7").Value  = 0.1143200680958068

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value  = 0.1935656666666667
$ws.Range("H8").Value  = 0.580697
$ws.Range("I8").Value  = 0.5985357597846208
$ws.Range("J8").Value  = 0.5985357597846208
$ws.Range("M8").Value  = 61.04160633333334
$ws.Range("N8").Value  = 183.124819
$ws.Range("O8").Value  = 0.2043613460574534
$ws.Range("P8").Value  = 0.2043613460574534
$ws.Range("Q8").Value  = 11.81555922431589
$ws.Range("R8").Value  = 106.340033018843
$ws.Range("S8").Value  = 0.1223175735331057
$ws.Range("T8").Value  = 0.1223175735331057

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value  = 0.1935656666666667
$ws.Range("H9").Value  = 0.580697
$ws.Range("I9").Value  = 0.5985357597846208
$ws.Range("J9").Value  = 0.5985357597846208
$ws.Range("O9").Value  = 0.3559304658284363
$ws.Range("P9").Value  = 0.3559304658284363
$ws.Range("Q9").Value  = 20.57883048760067
$ws.Range("R9").Value  = 185.209474388406
$ws.Range("S9").Value  = 0.2130371117951171
$ws.Range("T9").Value  = 0.2130371117951172

# Row 10 (MuSCs -> MuSCs)
$ws.Range("G10").Value = 0.1935656666666667
$ws.Range("H10").Value = 0.580697
$ws.Range("I10").Value = 0.5985357597846208
$ws.Range("J10").Value = 0.5985357597846208
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 25.42260676154622
$ws.Range("R10").Value = 228.803460853916
$ws.Range("S10").Value = 0.2631810744563979
$ws.Range("T10").Value = 0.263181074456398

Write-Output "Updated LR-pair TPM values (rows 2-10)"
